$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price (D) column ranges that receive new
# numeric-looking values, so Excel stores them as text (matching the
# inlineStr/shared-string representation in the target workbook) instead
# of auto-converting them to numbers. Rows 4 and 46 are excluded since
# their D column is not modified by this update.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D45").NumberFormat = "@"
$ws.Range("D47:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '30.881.86'
$ws.Range("E2").Value = '  +0.50%  '

# Row 3
$ws.Range("D3").Value = '1.906.57'
$ws.Range("E3").Value = '  +0.68%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").Value = '238.87'
$ws.Range("E5").Value = '  -3.17%  '

# Row 6
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.04%  '

# Row 7
$ws.Range("D7").Value = '0.4904'
$ws.Range("E7").Value = '  -0.48%  '

# Row 8
$ws.Range("D8").Value = '0.2963'
$ws.Range("E8").Value = '  +0.19%  '

# Row 9
$ws.Range("D9").Value = '0.06758'
$ws.Range("E9").Value = '  -0.72%  '

# Row 10
$ws.Range("D10").Value = '1.918.47'
$ws.Range("E10").Value = '  +1.25%  '

# Row 11
$ws.Range("D11").Value = '17.05'
$ws.Range("E11").Value = '  -1.42%  '

# Row 12
$ws.Range("D12").Value = '0.07290'
$ws.Range("E12").Value = '  +0.46%  '

# Row 13
$ws.Range("D13").Value = '89.81'
$ws.Range("E13").Value = '  -2.61%  '

# Row 14
$ws.Range("D14").Value = '5.120'
$ws.Range("E14").Value = '  +0.49%  '

# Row 15
$ws.Range("D15").Value = '0.6699'
$ws.Range("E15").Value = '  -1.84%  '

# Row 16
$ws.Range("D16").Value = '30.861.15'
$ws.Range("E16").Value = '  +0.49%  '

# Row 17
$ws.Range("D17").Value = '0.000007944'
$ws.Range("E17").Value = '  -0.57%  '

# Row 18
$ws.Range("D18").Value = '13.48'
$ws.Range("E18").Value = '  +1.28%  '

# Row 19
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.02%  '

# Row 20
$ws.Range("D20").Value = '2.154.51'
$ws.Range("E20").Value = '  +0.66%  '

# Row 21
$ws.Range("D21").Value = '0.9995'
$ws.Range("E21").Value = '  -0.16%  '

# Row 22
$ws.Range("D22").Value = '5.103'
$ws.Range("E22").Value = '  +5.24%  '

# Row 23
$ws.Range("D23").Value = '207.25'
$ws.Range("E23").Value = '  +8.34%  '

# Row 24
$ws.Range("D24").Value = '6.200'
$ws.Range("E24").Value = '  +2.27%  '

# Row 25
$ws.Range("D25").Value = '9.641'
$ws.Range("E25").Value = '  +2.54%  '

# Row 26
$ws.Range("D26").Value = '157.92'
$ws.Range("E26").Value = '  +1.30%  '

# Row 27
$ws.Range("D27").Value = '18.89'
$ws.Range("E27").Value = '  -1.35%  '

# Row 28
$ws.Range("D28").Value = '1.967'
$ws.Range("E28").Value = '  +2.02%  '

# Row 29
$ws.Range("D29").Value = '1.431'
$ws.Range("E29").Value = '  +2.24%  '

# Row 30
$ws.Range("D30").Value = '4.326'
$ws.Range("E30").Value = '  -1.20%  '

# Row 31
$ws.Range("D31").Value = '0.09169'
$ws.Range("E31").Value = '  +1.71%  '

# Row 32
$ws.Range("D32").Value = '4.050'
$ws.Range("E32").Value = '  +0.11%  '

# Row 33
$ws.Range("D33").Value = '0.05176'
$ws.Range("E33").Value = '  -0.60%  '

# Row 34
$ws.Range("D34").Value = '0.7521'
$ws.Range("E34").Value = '  +0.74%  '

# Row 35
$ws.Range("D35").Value = '1.114'
$ws.Range("E35").Value = '  -1.07%  '

# Row 36
$ws.Range("D36").Value = '2.710'
$ws.Range("E36").Value = '  -0.52%  '

# Row 37
$ws.Range("D37").Value = '0.01838'
$ws.Range("E37").Value = '  -1.30%  '

# Row 38
$ws.Range("D38").Value = '2.722'
$ws.Range("E38").Value = '  +1.73%  '

# Row 39
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '0.9271'
$ws.Range("E39").Value = '  -1.22%  '

# Row 40
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '2.102'
$ws.Range("E40").Value = '  -2.88%  '

# Row 41
$ws.Range("D41").Value = '0.4488'
$ws.Range("E41").Value = '  +1.02%  '

# Row 42
$ws.Range("D42").Value = '106.29'
$ws.Range("E42").Value = '  -0.18%  '

# Row 43
$ws.Range("D43").Value = '5.840'
$ws.Range("E43").Value = '  +0.64%  '

# Row 44
$ws.Range("D44").Value = '1.006'
$ws.Range("E44").Value = '  +0.60%  '

# Row 45
$ws.Range("D45").Value = '7.762'
$ws.Range("E45").Value = '  +0.85%  '

# Row 46
$ws.Range("E46").Value = '  +2.36%  '

# Row 47
$ws.Range("D47").Value = '67.18'
$ws.Range("E47").Value = '  +15.65%  '

# Row 48
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = '0.4109'
$ws.Range("E48").Value = '  +3.82%  '

# Row 49
$ws.Range("D49").Value = '35.01'
$ws.Range("E49").Value = '  +4.17%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '8.931'
$ws.Range("E50").Value = '  +2.06%  '

# Row 51
$ws.Range("D51").Value = '0.05905'
$ws.Range("E51").Value = '  +0.83%  '
